# Front-end Team Meeting #4 - apply "Add files via upload" edit
#
# Summary of change:
#  - Slide 1 ("Front-end Team Meeting #4"): merge the "25" / "/04/2019" runs
#    in the subtitle into a single run "25/04/2019".
#  - Slide 3 ("React-js Project 1: Todo List"): duplicated to become the new
#    slide 4 (with small text tweaks), and slide 3 itself is replaced with a
#    new "Reference links" slide (two hyperlinks).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Slide 1: merge date runs "25" + "/04/2019" -> "25/04/2019"
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$dateShape = $slide1.Shapes.Item(2)
# Force a real text replacement so the two existing runs collapse into one.
$dateShape.TextFrame.TextRange.Text = "x"
$dateShape.TextFrame.TextRange.Text = "25/04/2019"

# ---------------------------------------------------------------------
# 2. Duplicate slide 3 ("React-js Project 1: Todo List") so the original
#    content is preserved (with tweaks) on the new slide 4, freeing up
#    slide 3 to receive the new "Reference links" content.
# ---------------------------------------------------------------------
$slide3Old = $p.Slides.Item(3)
$dupRange = $slide3Old.Duplicate()
$slide4 = $p.Slides.Item(4)

# --- Slide 4: title -------------------------------------------------
$title4 = $slide4.Shapes.Item(1)
$title4.TextFrame.TextRange.Text = "x"
$title4.TextFrame.TextRange.Text = "React-js Project 1: Todo List"
$title4.TextFrame.TextRange.Font.Bold = -1

# --- Slide 4: content placeholder -------------------------------------
$content4 = $slide4.Shapes.Item(2)
$content4.Left = 36
$content4.Top = 108
$content4.Width = 648
$content4.Height = 374.3750787401575

$content4.TextFrame.TextRange.Text = "x"
$content4.TextFrame.TextRange.Text = "Complete this in 3 weeks.`rSubmit this project to git-hub before the next meeting 16 May 2019x`rGit-hub repository link:  https://github.com/tnnha76/reactjs-practice-projects"

# Paragraph 2: split trailing "x" placeholder into its own run, then turn
# it into the final period "." run (so "...2019" and "." are two runs).
$p2 = $content4.TextFrame.TextRange.Paragraphs(2, 1)
$p2len = $p2.Text.Length
$lastChar = $p2.Characters($p2len, 1)
$lastChar.Text = "."

# Paragraph 3: split "Git-hub " from "repository link:  " into two runs,
# and give the url run its hyperlink.
$p3 = $content4.TextFrame.TextRange.Paragraphs(3, 1)
$gitHubWord = $p3.Characters(1, 8)
$gitHubWord.Font.Bold = $gitHubWord.Font.Bold
$urlRun = $p3.Characters($p3.Text.Length - "https://github.com/tnnha76/reactjs-practice-projects".Length + 1, "https://github.com/tnnha76/reactjs-practice-projects".Length)
$urlRun.ActionSettings(1).Hyperlink.Address = "https://github.com/tnnha76/reactjs-practice-projects"

# ---------------------------------------------------------------------
# 3. Slide 3: replace with the new "Reference links" content.
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)

# --- Slide 3: title (unchanged text, just re-set cleanly) -----------
$title3 = $slide3.Shapes.Item(1)
$title3.TextFrame.TextRange.Text = "x"
$title3.TextFrame.TextRange.Text = "React-js Project 1: Todo List"
$title3.TextFrame.TextRange.Font.Bold = -1

# --- Slide 3: content placeholder ------------------------------------
$content3 = $slide3.Shapes.Item(2)
$content3.Left = 36
$content3.Top = 108
$content3.Width = 648
$content3.Height = 374.3750787401575

$viblo = "https://viblo.asia/p/hoc-reactjs-thong-qua-vi-du-phan-1-yMnKMnGaZ7P"
$youtube = "https://www.youtube.com/watch?v=494HPoUmqdI&list=PLJ5qtRQovuEOoKffoCBzTfvzMTTORnoyp&index=16"

$content3.TextFrame.TextRange.Text = "x"
$content3.TextFrame.TextRange.Text = "Reference links:`r$viblo `r$youtube`r"

# Paragraph 1: "Reference links:" has no bullet.
$content3.TextFrame.TextRange.Paragraphs(1, 1).ParagraphFormat.Bullet.Visible = 0

# Paragraph 2: hyperlink the viblo link (leave the trailing space plain).
$para2 = $content3.TextFrame.TextRange.Paragraphs(2, 1)
$vibloRun = $para2.Characters(1, $viblo.Length)
$vibloRun.ActionSettings(1).Hyperlink.Address = $viblo

# Paragraph 3: hyperlink the whole youtube link.
$para3 = $content3.TextFrame.TextRange.Paragraphs(3, 1)
$youtubeRun = $para3.Characters(1, $youtube.Length)
$youtubeRun.ActionSettings(1).Hyperlink.Address = $youtube

# Paragraph 4: trailing empty paragraph, also with no bullet.
$content3.TextFrame.TextRange.Paragraphs(4, 1).ParagraphFormat.Bullet.Visible = 0
